$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data block (EM / padding exercise), columns AG:AK, rows 5-9.
$values = @{
    "AG5" = 48;  "AH5" = 65;  "AI5" = 57;  "AJ5" = 63;  "AK5" = 90;
    "AG6" = 56;  "AH6" = 73;  "AI6" = 61;  "AJ6" = 60;  "AK6" = 73;
    "AG7" = 46;  "AH7" = 56;  "AI7" = 57;  "AJ7" = 61;  "AK7" = 83;
    "AG8" = 54;  "AH8" = 55;  "AI8" = 58;  "AJ8" = 80;  "AK8" = 89;
    "AG9" = 45;  "AH9" = 68;  "AI9" = 73;  "AJ9" = 74;  "AK9" = 86;
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Normal (black) bold Courier 12 style for the new block, except AI7.
$normalRange = $ws.Range("AG5:AK9")
$normalRange.HorizontalAlignment = -4108
$normalRange.Font.Size = 12
$normalRange.Font.Bold = $true

# AI7 is highlighted in red to call out the example value.
$ws.Range("AI7").Font.Color = 255

# Match the author's final selection from the commit.
$ws.Range("AJ9:AK9").Select()
